$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.899.43'
$ws.Range("E2").Value = '  +0.41%  '

# Row 3
$ws.Range("D3").Value = '2.260.21'
$ws.Range("E3").Value = '  -0.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.30'
$ws.Range("E5").Value = '  +0.93%  '

# Row 6
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.653'
$ws.Range("E6").Value = '  +4.24%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.48'
$ws.Range("E7").Value = '  -0.76%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.450'
$ws.Range("E9").Value = '  +5.59%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0975'
$ws.Range("E10").Value = '  -6.52%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.03'
$ws.Range("E11").Value = '  +1.27%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.48'
$ws.Range("E12").Value = '  +1.82%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("E13").Value = '  +1.73%  '

# Row 14
$ws.Range("D14").Value = '2.602.45'
$ws.Range("E14").Value = '  -0.28%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.60'
$ws.Range("E15").Value = '  -0.35%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.14'
$ws.Range("E16").Value = '  +4.45%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.841'
$ws.Range("E17").Value = '  +2.86%  '

# Row 18
$ws.Range("D18").Value = '2.259.72'
$ws.Range("E18").Value = '  -0.45%  '

# Row 19
$ws.Range("D19").Value = '43.864.00'
$ws.Range("E19").Value = '  +0.38%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0983'
$ws.Range("E20").Value = '  -1.67%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.67'
$ws.Range("E21").Value = '  +0.65%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").Value = '  +1.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.35'
$ws.Range("E23").Value = '  +0.25%  '

# Row 24
$ws.Range("E24").Value = '  -0.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("E25").Value = '  -1.35%  '

# Row 26
$ws.Range("B26").Value = 'WEMIXToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.45'
$ws.Range("E26").Value = '  +23.70%  '

# Row 27
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.31'
$ws.Range("E27").Value = '  -2.67%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.88'
$ws.Range("E28").Value = '  +0.28%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.95'
$ws.Range("E29").Value = '  +1.37%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.91'
$ws.Range("E30").Value = '  +4.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.44'
$ws.Range("E32").Value = '  -0.69%  '

# Row 33
$ws.Range("E33").Value = '  +2.76%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.95'
$ws.Range("E34").Value = '  +5.43%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0686'
$ws.Range("E35").Value = '  -0.42%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.94'
$ws.Range("E36").Value = '  -3.24%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.70'
$ws.Range("E37").Value = '  -2.96%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.41'
$ws.Range("E38").Value = '  -4.90%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.29'
$ws.Range("E39").Value = '  -1.58%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0255'
$ws.Range("E40").Value = '  +3.12%  '

# Row 41
$ws.Range("E41").Value = '  -0.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.69'
$ws.Range("E42").Value = '  +3.52%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.23'
$ws.Range("E43").Value = '  +0.00%  '

# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.57'
$ws.Range("E44").Value = '  +1.31%  '

# Row 45
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.40'
$ws.Range("E45").Value = '  -0.41%  '

# Row 46
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0950'
$ws.Range("E46").Value = '  -1.32%  '

# Row 47
$ws.Range("B47").Value = 'TerraClassic'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000212'
$ws.Range("E47").Value = '  +1.78%  '

# Row 48
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.19'
$ws.Range("E48").Value = '  -1.05%  '

# Row 49
$ws.Range("D49").Value = '1.454.41'
$ws.Range("E49").Value = '  -1.47%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.33'
$ws.Range("E50").Value = '  +0.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.91'
$ws.Range("E51").Value = '  -6.69%  '
